$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.211.69'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.587.97'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.19%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.87'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.86%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.505'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.20%  '
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0604'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.14%  '
$ws.Range("E10").Value = '  -1.69%  '
$ws.Range("E11").Value = '  +0.23%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.813.19'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.13%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.586.08'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.77%  '
$ws.Range("E14").Value = '  -1.44%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.98'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.01%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.224.89'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.31%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0₃0725'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.56%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '214.44'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.36%  '
$ws.Range("E20").Value = '  -2.76%  '
$ws.Range("E21").Value = '  -0.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.24'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.83%  '
$ws.Range("E23").Value = '  -0.72%  '
$ws.Range("E24").Value = '  +0.42%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.77'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.18%  '
$ws.Range("E27").Value = '  -0.75%  '
$ws.Range("E28").Value = '  -0.73%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.14'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.91%  '
$ws.Range("E30").Value = '  -1.76%  '
$ws.Range("E31").Value = '  +0.80%  '
$ws.Range("E32").Value = '  -1.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.388.87'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.80%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.94'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.60%  '
$ws.Range("E35").Value = '  -0.51%  '
$ws.Range("E36").Value = '  -1.36%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.584'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.22%  '
$ws.Range("E38").Value = '  -0.71%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.818'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.71%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.85'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.06%  '
$ws.Range("E41").Value = '  -0.15%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.939'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -14.98%  '
$ws.Range("E43").Value = '  +0.62%  '
$ws.Range("E44").Value = '  -0.16%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.724.86'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.18%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '60.98'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.66%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '86.15'
$ws.Range("D47").Style = "Normal"
$ws.Range("E48").Value = '  -1.50%  '
$ws.Range("E49").Value = '  -0.87%  '
$ws.Range("E50").Value = '  -1.00%  '
$ws.Range("E51").Value = '  -0.23%  '
